# Update event stats on the "展览" and "全部类型" sheets.
# For each of these sheets:
#   F2: 1721 -> 1722
#   G2: 60 (number) -> "不可售" (text)
#   F4: 488 -> 489
#   F6: 80  -> 81
#   F7: 663 -> 669
#   F8: 411 -> 412

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1722
    $ws.Range("G2").Value = "不可售"

    $ws.Range("F4").Value = 489
    $ws.Range("F6").Value = 81
    $ws.Range("F7").Value = 669
    $ws.Range("F8").Value = 412
}
